$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: "...inspect the site on a XXXX basis..."
#   XXXX (bold, red)  ->  { VISITS } visits per month
#   ({, VISITS, }, " " stay bold+red; "visits per month" is bold, no red)
# ---------------------------------------------------------------------
$r1 = $d.Content
$null = $r1.Find.Execute("XXXX basis", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start1 = $r1.Start

# Replace just the "XXXX" token (single run) with the full literal text;
# this keeps the original run's bold+red formatting for everything.
$xxxx = $d.Range($start1, $start1 + 4)
$xxxx.Text = "{VISITS} visits per month"

# Now strip the red colour from the trailing "visits per month" part only,
# matching the plain (non-red) run used by the rest of the sentence.
$tail = $d.Range($start1 + 9, $start1 + 9 + 16)
$tail.Font.Color = -16777216

# ---------------------------------------------------------------------
# Edit 2: "...Construction Permit fee of RX X XXX-XX excl. VAT..."
#   "X" + "X XXX-XX" (identically formatted, red) -> {PERMIT_FEE}
# ---------------------------------------------------------------------
$r2 = $d.Content
$null = $r2.Find.Execute("XX XXX-XX", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2.Text = "{PERMIT_FEE}"

# ---------------------------------------------------------------------
# Edit 3: "...the duration of which will be XX (XX) months..."
#   Only the second "XX" (inside the parentheses) -> {TOTAL_MONTHS}
# ---------------------------------------------------------------------
$r3 = $d.Content
$null = $r3.Find.Execute("(XX) months", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start3 = $r3.Start
$inner3 = $d.Range($start3 + 1, $start3 + 3)
$inner3.Text = "{TOTAL_MONTHS}"

# ---------------------------------------------------------------------
# Edit 4: "Upon acceptance of this proposal (XX-XXXX-XX-XX22) - ..."
#   "XX" + "-XXXX-" + "XX" + "-XX2" + "2" (identically formatted, red/bold/italic)
#   -> {REF}
# ---------------------------------------------------------------------
$r4 = $d.Content
$null = $r4.Find.Execute("(XX-XXXX-XX-XX22)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start4 = $r4.Start
$end4 = $r4.End
$inner4 = $d.Range($start4 + 1, $end4 - 1)
$inner4.Text = "{REF}"
